$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.625513630481521
$ws.Range("C2").Value = 0.6617815912354104
$ws.Range("D2").Value = 0.1179493099002187
$ws.Range("F2").Value = 2.885448188488681
$ws.Range("G2").Value = 0.002511245253503473
$ws.Range("I2").Value = 1.263084568252836
$ws.Range("J2").Value = 0.3530977822467776

$ws.Range("B3").Value = 1.505060230384004
$ws.Range("C3").Value = 0.6112078090272348
$ws.Range("D3").Value = 0.1165508021603472
$ws.Range("F3").Value = 2.834541716054702
$ws.Range("G3").Value = 0.002517378919973599
$ws.Range("I3").Value = 1.250980048515537
$ws.Range("J3").Value = 0.3423358023735545

$ws.Range("B4").Value = 1.431977383770743
$ws.Range("C4").Value = 0.5805467742078463
$ws.Range("D4").Value = 0.1157374678381373
$ws.Range("F4").Value = 2.805232752602635
$ws.Range("G4").Value = 0.002521338183833507
$ws.Range("I4").Value = 1.244390930733239
$ws.Range("J4").Value = 0.3359622724288869

$ws.Range("B5").Value = 1.40241442024427
$ws.Range("C5").Value = 0.5681493427827604
$ws.Range("D5").Value = 0.1154175136749913
$ws.Range("F5").Value = 2.793774736636792
$ws.Range("G5").Value = 0.002523000371929517
$ws.Range("I5").Value = 1.24191590869254
$ws.Range("J5").Value = 0.3334234734701766

$ws.Range("B6").Value = 1.397518693772724
$ws.Range("C6").Value = 0.5660965927917232
$ws.Range("D6").Value = 0.1153650816904275
$ws.Range("F6").Value = 2.791901351087503
$ws.Range("G6").Value = 0.002523279326780732
$ws.Range("I6").Value = 1.241517566852146
$ws.Range("J6").Value = 0.3330054255339121

$ws.Range("B7").Value = 1.431577802817458
$ws.Range("C7").Value = 0.5803791862261392
$ws.Range("D7").Value = 0.1157331062068678
$ws.Range("F7").Value = 2.805076264989623
$ws.Range("G7").Value = 0.002521360403180684
$ws.Range("I7").Value = 1.244356703556292
$ws.Range("J7").Value = 0.335927797182265

$ws.Range("B8").Value = 1.583798578820506
$ws.Range("C8").Value = 0.6442617244570101
$ws.Range("D8").Value = 0.117457738473469
$ws.Range("F8").Value = 2.867488826413137
$ws.Range("G8").Value = 0.002513320159405392
$ws.Range("I8").Value = 1.258734785310423
$ws.Range("J8").Value = 0.3493381204471859

$ws.Range("B9").Value = 1.889336880257417
$ws.Range("C9").Value = 0.7727057202478704
$ws.Range("D9").Value = 0.1211964026996171
$ws.Range("F9").Value = 3.00553315596656
$ws.Range("G9").Value = 0.002499077640865894
$ws.Range("I9").Value = 1.293708821115061
$ws.Range("J9").Value = 0.3775184696080771

$ws.Range("B10").Value = 2.118247715161601
$ws.Range("C10").Value = 0.8691071218896695
$ws.Range("D10").Value = 0.1241568587028041
$ws.Range("F10").Value = 3.116786938427794
$ws.Range("G10").Value = 0.002489531263280627
$ws.Range("I10").Value = 1.323664124440938
$ws.Range("J10").Value = 0.3994053242518874

$ws.Range("B11").Value = 2.223381514409311
$ws.Range("C11").Value = 0.9134271459662386
$ws.Range("D11").Value = 0.1255493069411813
$ws.Range("F11").Value = 3.169598385337423
$ws.Range("G11").Value = 0.002485385108299987
$ws.Range("I11").Value = 1.338244526890676
$ws.Range("J11").Value = 0.4096269157389969

$ws.Range("B12").Value = 2.263339252122194
$ws.Range("C12").Value = 0.9302787769729548
$ws.Range("D12").Value = 0.1260830989806863
$ws.Range("F12").Value = 3.189918312235903
$ws.Range("G12").Value = 0.002483843136840496
$ws.Range("I12").Value = 1.343905096993979
$ws.Range("J12").Value = 0.4135362974995331

$ws.Range("B13").Value = 2.254727120616053
$ws.Range("C13").Value = 0.9266464012306415
$ws.Range("D13").Value = 0.1259678492651659
$ws.Range("F13").Value = 3.185527681890164
$ws.Range("G13").Value = 0.002484173981562824
$ws.Range("I13").Value = 1.342679762354493
$ws.Range("J13").Value = 0.4126926125618695

$ws.Range("B14").Value = 2.226665927350268
$ws.Range("C14").Value = 0.9148121561469225
$ws.Range("D14").Value = 0.1255930924296678
$ws.Range("F14").Value = 3.171263649332872
$ws.Range("G14").Value = 0.002485257687580992
$ws.Range("I14").Value = 1.338707421382864
$ws.Range("J14").Value = 0.4099477643445368

$ws.Range("B15").Value = 2.209496687775754
$ws.Range("C15").Value = 0.9075723185957827
$ws.Range("D15").Value = 0.1253643876944466
$ws.Range("F15").Value = 3.162568514002658
$ws.Range("G15").Value = 0.002485925141513105
$ws.Range("I15").Value = 1.336292451875906
$ws.Range("J15").Value = 0.4082715205104677

$ws.Range("B16").Value = 2.111397278707727
$ws.Range("C16").Value = 0.8662202204029086
$ws.Range("D16").Value = 0.1240667727189049
$ws.Range("F16").Value = 3.113380295558727
$ws.Range("G16").Value = 0.002489806164075562
$ws.Range("I16").Value = 1.322730633858015
$ws.Range("J16").Value = 0.3987427022037906

$ws.Range("B17").Value = 2.051474195288392
$ws.Range("C17").Value = 0.8409726030615161
$ws.Range("D17").Value = 0.1232823860828915
$ws.Range("F17").Value = 3.083772178395407
$ws.Range("G17").Value = 0.002492237259890859
$ws.Range("I17").Value = 1.314656666058696
$ws.Range("J17").Value = 0.3929653860937634

$ws.Range("B18").Value = 2.017102133566368
$ws.Range("C18").Value = 0.8264946617181863
$ws.Range("D18").Value = 0.1228355354122925
$ws.Range("F18").Value = 3.06694930132025
$ws.Range("G18").Value = 0.002493654071539357
$ws.Range("I18").Value = 1.310102358522727
$ws.Range("J18").Value = 0.3896673498827283

$ws.Range("B19").Value = 2.005480451518622
$ws.Range("C19").Value = 0.8216001603477707
$ws.Range("D19").Value = 0.1226849815650581
$ws.Range("F19").Value = 3.061288759250374
$ws.Range("G19").Value = 0.002494136963914589
$ws.Range("I19").Value = 1.308575675884015
$ws.Range("J19").Value = 0.3885549552538521

$ws.Range("B20").Value = 2.057843356270894
$ws.Range("C20").Value = 0.8436557089663097
$ws.Range("D20").Value = 0.1233654400735844
$ws.Range("F20").Value = 3.086902562121395
$ws.Range("G20").Value = 0.00249197655114218
$ws.Range("I20").Value = 1.315506862249123
$ws.Range("J20").Value = 0.3935778077916439

$ws.Range("B21").Value = 2.234904206201122
$ws.Range("C21").Value = 0.9182862867697281
$ws.Range("D21").Value = 0.1257029916995123
$ws.Range("F21").Value = 3.175444580771057
$ws.Range("G21").Value = 0.002484938615990209
$ws.Range("I21").Value = 1.339870396616007
$ws.Range("J21").Value = 0.4107529383414459

$ws.Range("B22").Value = 2.351474770998493
$ws.Range("C22").Value = 0.9674622661668764
$ws.Range("D22").Value = 0.127268583894363
$ws.Range("F22").Value = 3.235187421201772
$ws.Range("G22").Value = 0.002480502550635765
$ws.Range("I22").Value = 1.356606273469424
$ws.Range("J22").Value = 0.4222036436159442

$ws.Range("B23").Value = 2.289180346573062
$ws.Range("C23").Value = 0.9411789571068994
$ws.Range("D23").Value = 0.126429556487281
$ws.Range("F23").Value = 3.203128330199036
$ws.Range("G23").Value = 0.002482855247574529
$ws.Range("I23").Value = 1.347598919684017
$ws.Range("J23").Value = 0.4160713416051891

$ws.Range("B24").Value = 2.054963614168855
$ws.Range("C24").Value = 0.8424425608752699
$ws.Range("D24").Value = 0.1233278785769443
$ws.Range("F24").Value = 3.085486695463146
$ws.Range("G24").Value = 0.002492094357847829
$ws.Range("I24").Value = 1.315122216096881
$ws.Range("J24").Value = 0.393300859036728

$ws.Range("B25").Value = 1.80591309161764
$ws.Range("C25").Value = 0.7376085110415715
$ws.Range("D25").Value = 0.1201472094188674
$ws.Range("F25").Value = 2.966482947265121
$ws.Range("G25").Value = 0.002502768635538069
$ws.Range("I25").Value = 1.263508615101923
$ws.Range("J25").Value = 0.3696898009208667
